# Add season-record columns (Wins, Losses, Ties) to the player table.
# The previous export only captured team statistics, not the season
# win/loss/tie record, so three new columns are appended after the
# existing data (column AB) for every row of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 45

# New header cells, appended right after the existing "Unnamed: 27" (AB1) header.
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the formatting of the rest of the header row (bold, centered, bordered)
# by copying the style from the neighboring header cell instead of the default.
$headerStyleSource = $ws.Range("AB1")
$newHeaderCells = $ws.Range("AC1:AE1")
$headerStyleSource.Copy()
$newHeaderCells.PasteSpecial(-4122)

# Season record values for this team: 64 wins, 98 losses, 0 ties.
$wins = 64
$losses = 98
$ties = 0

for ($row = 2; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, 29).Value = $wins
    $ws.Cells.Item($row, 30).Value = $losses
    $ws.Cells.Item($row, 31).Value = $ties
}
